$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to remain plain text so values like
# "7.90", "0.910", "3.50" keep their trailing zeros instead of being
# auto-coerced to numbers by Excels type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "95.360.85"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "3.595.07"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "235.01"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "653.48"
$ws.Range("E6").Value = "  +4.88%  "
$ws.Range("D7").Value = "1.45"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").Value = "0.399"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "3.593.51"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "42.14"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").Value = "4.286.10"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Value = "95.180.75"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "3.598.96"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "7.90"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "17.87"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "3.50"
$ws.Range("E22").Value = "  +3.48%  "
$ws.Range("D23").Value = "506.37"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "0.473"
$ws.Range("E24").Value = "  -4.96%  "
$ws.Range("E25").Value = "  +6.06%  "
$ws.Range("D26").Value = "6.57"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").Value = "91.95"
$ws.Range("E27").Value = "  -4.29%  "
$ws.Range("D28").Value = "3.788.41"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").Value = "12.45"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").Value = "3.04"
$ws.Range("E30").Value = "  +4.58%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "11.20"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "0.138"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "32.42"
$ws.Range("E35").Value = "  +8.97%  "
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").Value = "0.559"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "8.05"
$ws.Range("E38").Value = "  +6.89%  "
$ws.Range("D39").Value = "555.97"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "0.910"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "35.88"
$ws.Range("E44").Value = "  +40.05%  "
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "2.28"
$ws.Range("E46").Value = "  +6.50%  "
$ws.Range("D47").Value = "23.55"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "5.65"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("D49").Value = "0.0411"
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("E51").Value = "  +0.07%  "
